$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.3
$ws.Range("J2").Value = 3
$ws.Range("P2").Value = 1.64
$ws.Range("Q2").Value = 2.24

# Row 3
$ws.Range("F3").Value = 1.63
$ws.Range("G3").Value = 1.69
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 7.4
$ws.Range("K3").Value = 4.1

# Row 4
$ws.Range("P4").Value = 1.25

# Row 5
$ws.Range("P5").Value = 1.25
